# Update the "dSF" (column F) values for specific rows, per the repull/recalculation
# of the underlying data described in the commit message.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = -12
$ws.Range("F4").Value = 6
$ws.Range("F6").Value = -13
$ws.Range("F7").Value = -5
$ws.Range("F8").Value = -5
$ws.Range("F10").Value = -1
$ws.Range("F13").Value = -2
$ws.Range("F25").Value = 1
